$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): columns C, D, E get new labels
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2-8: column C becomes the family text (like D), column E becomes numeric 1
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = "f__Anaeroplasmataceae"
    $ws.Cells.Item($r, 5).Value = 1
}
